$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 -> becomes the "JOSE GABRIEL AVILA PADILLA" / period 2308 record
$ws.Range("C16").Value = "1062675657"
$ws.Range("D16").Value = "JOSE GABRIEL AVILA PADILLA"
$ws.Range("E16").Value = "2308"
$ws.Range("F16").Value = 82388
$ws.Range("G16").Value = 2308400

# Row 18 -> becomes the "CARLOS ANDRES ESPITIA ROMERO" / period 2303 record
$ws.Range("C18").Value = "1067918429"
$ws.Range("D18").Value = "CARLOS ANDRES ESPITIA ROMERO"
$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 105611
$ws.Range("G18").Value = 4227100
